# Updates the "cryptos" price/volume table (GitHub Actions refresh snapshot):
#  - refreshes Price (D) / Volume(1h) (E) text for most rows
#  - Filecoin <-> SuiNetwork (rows 41/42) and Maker <-> RenderToken (rows 49/50)
#    swapped rank order, so their Coin/Link/Price/Volume cells trade places
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on price cells whose new value would otherwise be
# auto-converted to a number by Excel (matches the source inlineStr text).
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"

# Apply the updated cell values
$ws.Range("D2").Value = '60.764.89'
$ws.Range("E2").Value = '  +6.45%  '
$ws.Range("D3").Value = '2.668.12'
$ws.Range("E3").Value = '  +10.72%  '
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("D5").Value = '513.94'
$ws.Range("E5").Value = '  +5.26%  '
$ws.Range("D6").Value = '158.83'
$ws.Range("E6").Value = '  +2.59%  '
$ws.Range("D7").Value = '0.998'
$ws.Range("E7").Value = '  +0.21%  '
$ws.Range("D8").Value = '0.605'
$ws.Range("E8").Value = '  -0.17%  '
$ws.Range("D9").Value = '2.693.31'
$ws.Range("E9").Value = '  +10.98%  '
$ws.Range("D10").Value = '6.47'
$ws.Range("E10").Value = '  +1.26%  '
$ws.Range("D11").Value = '0.106'
$ws.Range("E11").Value = '  +5.74%  '
$ws.Range("D12").Value = '0.349'
$ws.Range("E12").Value = '  +3.75%  '
$ws.Range("E13").Value = '  +1.19%  '
$ws.Range("D14").Value = '3.151.55'
$ws.Range("E14").Value = '  +11.46%  '
$ws.Range("D15").Value = '60.999.02'
$ws.Range("E15").Value = '  +6.96%  '
$ws.Range("D16").Value = '21.99'
$ws.Range("E16").Value = '  +5.86%  '
$ws.Range("D17").Value = '0.0000142'
$ws.Range("E17").Value = '  +5.99%  '
$ws.Range("D18").Value = '2.685.16'
$ws.Range("E18").Value = '  +10.89%  '
$ws.Range("E19").Value = '  +1.51%  '
$ws.Range("D20").Value = '350.85'
$ws.Range("E20").Value = '  +7.83%  '
$ws.Range("D21").Value = '10.60'
$ws.Range("E21").Value = '  +6.35%  '
$ws.Range("D22").Value = '6.24'
$ws.Range("E22").Value = '  +4.29%  '
$ws.Range("D23").Value = '1.00'
$ws.Range("E23").Value = '  +0.18%  '
$ws.Range("D24").Value = '60.63'
$ws.Range("E24").Value = '  +4.03%  '
$ws.Range("D25").Value = '0.427'
$ws.Range("E25").Value = '  +4.96%  '
$ws.Range("D26").Value = '2.802.89'
$ws.Range("E26").Value = '  +11.06%  '
$ws.Range("D27").Value = '0.168'
$ws.Range("E27").Value = '  +4.55%  '
$ws.Range("D28").Value = '0.997'
$ws.Range("E28").Value = '  +0.03%  '
$ws.Range("D29").Value = '0.0₃0876'
$ws.Range("E29").Value = '  +11.75%  '
$ws.Range("D30").Value = '7.59'
$ws.Range("E30").Value = '  +4.11%  '
$ws.Range("D32").Value = '19.65'
$ws.Range("E32").Value = '  +5.59%  '
$ws.Range("D33").Value = '157.28'
$ws.Range("E33").Value = '  +4.86%  '
$ws.Range("E34").Value = '  +3.85%  '
$ws.Range("D35").Value = '5.79'
$ws.Range("E35").Value = '  +9.50%  '
$ws.Range("D36").Value = '4.10'
$ws.Range("E36").Value = '  +9.88%  '
$ws.Range("D37").Value = '1.23'
$ws.Range("E37").Value = '  +5.58%  '
$ws.Range("D38").Value = '315.26'
$ws.Range("E38").Value = '  +17.00%  '
$ws.Range("D39").Value = '1.52'
$ws.Range("E39").Value = '  +10.32%  '
$ws.Range("D40").Value = '0.867'
$ws.Range("E40").Value = '  +2.65%  '
$ws.Range("B41").Value = 'SuiNetwork'
$ws.Range("C41").Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range("D41").Value = '0.847'
$ws.Range("E41").Value = '  +32.28%  '
$ws.Range("B42").Value = 'Filecoin'
$ws.Range("C42").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D42").Value = '3.79'
$ws.Range("E42").Value = '  +7.01%  '
$ws.Range("D43").Value = '35.46'
$ws.Range("E43").Value = '  +3.68%  '
$ws.Range("D44").Value = '0.648'
$ws.Range("E44").Value = '  +8.67%  '
$ws.Range("D45").Value = '0.0582'
$ws.Range("E45").Value = '  +9.23%  '
$ws.Range("E46").Value = '  -0.24%  '
$ws.Range("D47").Value = '20.23'
$ws.Range("E47").Value = '  +15.99%  '
$ws.Range("D48").Value = '0.993'
$ws.Range("E48").Value = '  -0.19%  '
$ws.Range("B49").Value = 'RenderToken'
$ws.Range("C49").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D49").Value = '4.95'
$ws.Range("E49").Value = '  +8.61%  '
$ws.Range("B50").Value = 'Maker'
$ws.Range("C50").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D50").Value = '2.082.89'
$ws.Range("E50").Value = '  +10.98%  '
$ws.Range("E51").Value = '  +3.30%  '
